# Update column G ("K") values for rows 2-20 on Sheet1, per the regen of
# save_data (switch Strike# -> K, recalculated values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 2
    4  = 6
    5  = 3
    6  = 4
    7  = 3
    8  = 5
    9  = 4
    10 = 3
    11 = 13
    12 = 0
    13 = 2
    14 = 6
    15 = 3
    16 = 3
    17 = 9
    18 = 2
    19 = 1
    20 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
